$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting A:E to B:F
$ws.Range("A1").EntireColumn.Insert()

# Copy header formatting (bold/centered/bordered) from B1 onto the new A1,
# then set its text
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Value = "ID"

# New ID values for rows 2-23
$ids = @(
    "H 4",
    "H 17",
    "H 66a",
    "H 72",
    "H 105",
    "H 106a",
    "H 154",
    "H 269",
    "H 270",
    "H 271",
    "H 274",
    "H 275",
    "H 1902 Grube 56 I-IV",
    "H 1904 III/1",
    "H 1912 XIII/3",
    "H 1912 XIV.3",
    "H 1936 DIV",
    "H 1968 Sk 4",
    "H 1968 Sk 6",
    "H 12301 K I/1",
    "H 12370 K IX/6",
    "H 12640 K XIII/3"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
